# Update "想去人数" (F column) counts on sheets "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 396
$ws1.Range("F9").Value = 6671
$ws1.Range("F12").Value = 147
$ws1.Range("F16").Value = 15967
$ws1.Range("F19").Value = 316
$ws1.Range("F22").Value = 11236
$ws1.Range("F24").Value = 4414

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 396
$ws4.Range("F10").Value = 6671
$ws4.Range("F13").Value = 147
$ws4.Range("F18").Value = 15967
$ws4.Range("F21").Value = 316
$ws4.Range("F25").Value = 11236
$ws4.Range("F27").Value = 4414
